# Update the worksheet per the commit: rename component names
# X_CH -> X_PG, X_LI -> X_TAG, S_F -> S_G (header row D1, E1, H1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "X_PG"
$ws.Range("E1").Value = "X_TAG"
$ws.Range("H1").Value = "S_G"

# Update selection to H1 as in the diff
$ws.Range("H1").Select()
